# Update the software BOM to reflect the currently used version of Unity.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Unity3D version (row 2, column C) was bumped from 2022.3.5f1 to 2022.3.13f1.
$ws.Range("C2").Value = "2022.3.13f1"

# Move the active selection to the cell that was just edited.
$ws.Range("C2").Select()
